$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.147.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.396.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.396.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.122"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.977.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.400.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.227.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.531.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.551"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000123"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.174"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.66%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.50%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0768"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.776"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.520.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -0.11%  "
